# refatorando o consolidador para modelo ETL
# Update the absenteeism dataset rows (2-11) with refreshed ETL output values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2727,  "Emanuella Gonçalves",   "TI",                     "Viagem de negócios", 4, 45092, 10077.26),
    @(77844, "Samuel Moreira",        "Recursos Humanos",       "Viagem de negócios", 4, 45081, 6297.91),
    @(38990, "Juliana das Neves",     "Atendimento ao Cliente", "Consulta médica",    6, 45102, 10889.64),
    @(91830, "Lorena Moura",          "P&D",                    "Viagem de negócios", 2, 45099, 11922.9),
    @(22631, "Alice Souza",           "Vendas",                 "Doença",             1, 45095, 8721.92),
    @(32392, "Giovanna da Cunha",     "Atendimento ao Cliente", "Viagem de negócios", 3, 45083, 2738.37),
    @(12967, "Brenda Ferreira",       "Financeiro",              "Viagem de negócios", 3, 45086, 6341.91),
    @(48251, "Joaquim Campos",        "Recursos Humanos",       "Problemas pessoais", 4, 45095, 12194.6),
    @(39378, "Luiz Otávio Nogueira",  "TI",                     "Consulta médica",    3, 45087, 8490.700000000001),
    @(92452, "Alexandre Duarte",      "Vendas",                 "Consulta médica",    7, 45092, 4987.31)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
